$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("G2").Value = 42.33019849999999
$ws.Range("H2").Value = 84.66039699999999
$ws.Range("I2").Value = 0.03880108177208133
$ws.Range("J2").Value = 0.02631398137952337
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 0.500046
$ws.Range("N2").Value = 1.000092
$ws.Range("O2").Value = 0.02009953378395692
$ws.Range("P2").Value = 0.01521121714331422
$ws.Range("Q2").Value = 21.167046439131
$ws.Range("R2").Value = 84.66818575652398
$ws.Range("S2").Value = 0.0007798836539320239
$ws.Range("T2").Value = 0.000400267684669057

$ws.Range("E3").Value = 2
$ws.Range("G3").Value = 42.33019849999999
$ws.Range("H3").Value = 84.66039699999999
$ws.Range("I3").Value = 0.03880108177208133
$ws.Range("J3").Value = 0.02631398137952337
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.821983
$ws.Range("N3").Value = 2.465949
$ws.Range("O3").Value = 0.03303991048491191
$ws.Range("P3").Value = 0.03750663509291002
$ws.Range("Q3").Value = 34.7947035536255
$ws.Range("R3").Value = 208.768221321753
$ws.Range("S3").Value = 0.001281984268467314
$ws.Range("T3").Value = 0.0009869488974434118

$ws.Range("E4").Value = 2
$ws.Range("G4").Value = 42.33019849999999
$ws.Range("H4").Value = 84.66039699999999
$ws.Range("I4").Value = 0.03880108177208133
$ws.Range("J4").Value = 0.02631398137952337
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 7.086962666666667
$ws.Range("N4").Value = 21.260888
$ws.Range("O4").Value = 0.2848630836849173
$ws.Range("P4").Value = 0.3233742335981926
$ws.Range("Q4").Value = 299.9925364420893
$ws.Range("R4").Value = 1799.955218652536
$ws.Range("S4").Value = 0.01105299580390572
$ws.Range("T4").Value = 0.008509263561520479

$ws.Range("E5").Value = 2
$ws.Range("G5").Value = 42.33019849999999
$ws.Range("H5").Value = 84.66039699999999
$ws.Range("I5").Value = 0.03880108177208133
$ws.Range("J5").Value = 0.02631398137952337
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 6.510773
$ws.Range("N5").Value = 19.532319
$ws.Range("O5").Value = 0.2617029270770581
$ws.Range("P5").Value = 0.2970830139841956
$ws.Range("Q5").Value = 275.6023134784405
$ws.Range("R5").Value = 1653.613880870643
$ws.Range("S5").Value = 0.01015435667350997
$ws.Range("T5").Value = 0.007817436898152803

$ws.Range("E6").Value = 2
$ws.Range("G6").Value = 42.33019849999999
$ws.Range("H6").Value = 84.66039699999999
$ws.Range("I6").Value = 0.03880108177208133
$ws.Range("J6").Value = 0.02631398137952337
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.570313333333333
$ws.Range("N6").Value = 4.71094
$ws.Range("O6").Value = 0.06311932481158002
$ws.Range("P6").Value = 0.07165253925551318
$ws.Range("Q6").Value = 66.47167510719666
$ws.Range("R6").Value = 398.83005064318
$ws.Range("S6").Value = 0.002449098083412678
$ws.Range("T6").Value = 0.001885463583765141

$ws.Range("E7").Value = 2
$ws.Range("G7").Value = 42.33019849999999
$ws.Range("H7").Value = 84.66039699999999
$ws.Range("I7").Value = 0.03880108177208133
$ws.Range("J7").Value = 0.02631398137952337
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 8.3884095
$ws.Range("N7").Value = 16.776819
$ws.Range("O7").Value = 0.3371752201575759
$ws.Range("P7").Value = 0.2551723609258745
$ws.Range("Q7").Value = 355.0830392342857
$ws.Range("R7").Value = 1420.332156937143
$ws.Range("S7").Value = 0.01308276328885363
$ws.Range("T7").Value = 0.006714600753972478

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 27.92162266666667
$ws.Range("H8").Value = 83.76486800000001
$ws.Range("I8").Value = 0.0255937652713472
$ws.Range("J8").Value = 0.02603563478222566
$ws.Range("K8").Value = 2
$ws.Range("M8").Value = 0.500046
$ws.Range("N8").Value = 1.000092
$ws.Range("O8").Value = 0.02009953378395692
$ws.Range("P8").Value = 0.01521121714331422
$ws.Range("Q8").Value = 13.962095727976
$ws.Range("R8").Value = 83.77257436785601
$ws.Range("S8").Value = 0.0005144227497301064
$ws.Range("T8").Value = 0.000396033694136459

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 27.92162266666667
$ws.Range("H9").Value = 83.76486800000001
$ws.Range("I9").Value = 0.0255937652713472
$ws.Range("J9").Value = 0.02603563478222566
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.821983
$ws.Range("N9").Value = 2.465949
$ws.Range("O9").Value = 0.03303991048491191
$ws.Range("P9").Value = 0.03750663509291002
$ws.Range("Q9").Value = 22.95109916441467
$ws.Range("R9").Value = 206.559892479732
$ws.Range("S9").Value = 0.0008456157135371587
$ws.Range("T9").Value = 0.0009765090531892136

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 27.92162266666667
$ws.Range("H10").Value = 83.76486800000001
$ws.Range("I10").Value = 0.0255937652713472
$ws.Range("J10").Value = 0.02603563478222566
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 7.086962666666667
$ws.Range("N10").Value = 21.260888
$ws.Range("O10").Value = 0.2848630836849173
$ws.Range("P10").Value = 0.3233742335981926
$ws.Range("Q10").Value = 197.8794974314205
$ws.Range("R10").Value = 1780.915476882784
$ws.Range("S10").Value = 0.007290718898303905
$ws.Range("T10").Value = 0.008419253443944669

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 27.92162266666667
$ws.Range("H11").Value = 83.76486800000001
$ws.Range("I11").Value = 0.0255937652713472
$ws.Range("J11").Value = 0.02603563478222566
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 6.510773
$ws.Range("N11").Value = 19.532319
$ws.Range("O11").Value = 0.2617029270770581
$ws.Range("P11").Value = 0.2970830139841956
$ws.Range("Q11").Value = 181.7913469743214
$ws.Range("R11").Value = 1636.122122768892
$ws.Range("S11").Value = 0.006697963286434717
$ws.Range("T11").Value = 0.007734744852095354

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 27.92162266666667
$ws.Range("H12").Value = 83.76486800000001
$ws.Range("I12").Value = 0.0255937652713472
$ws.Range("J12").Value = 0.02603563478222566
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 1.570313333333333
$ws.Range("N12").Value = 4.71094
$ws.Range("O12").Value = 0.06311932481158002
$ws.Range("P12").Value = 0.07165253925551318
$ws.Range("Q12").Value = 43.84569636176889
$ws.Range("R12").Value = 394.61126725592
$ws.Range("S12").Value = 0.0016154611833135
$ws.Range("T12").Value = 0.001865519343275629

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 27.92162266666667
$ws.Range("H13").Value = 83.76486800000001
$ws.Range("I13").Value = 0.0255937652713472
$ws.Range("J13").Value = 0.02603563478222566
$ws.Range("K13").Value = 2
$ws.Range("M13").Value = 8.3884095
$ws.Range("N13").Value = 16.776819
$ws.Range("O13").Value = 0.3371752201575759
$ws.Range("P13").Value = 0.2551723609258745
$ws.Range("Q13").Value = 234.218004832482
$ws.Range("R13").Value = 1405.308028994892
$ws.Range("S13").Value = 0.008629583440027811
$ws.Range("T13").Value = 0.006643574395584339

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 441.022868
$ws.Range("H14").Value = 1323.068604
$ws.Range("I14").Value = 0.4042542905775846
$ws.Range("J14").Value = 0.4112336327632385
$ws.Range("K14").Value = 2
$ws.Range("M14").Value = 0.500046
$ws.Range("N14").Value = 1.000092
$ws.Range("O14").Value = 0.02009953378395692
$ws.Range("P14").Value = 0.01521121714331422
$ws.Range("Q14").Value = 220.531721051928
$ws.Range("R14").Value = 1323.190326311568
$ws.Range("S14").Value = 0.008125322770773699
$ws.Range("T14").Value = 0.006255364084595559

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 441.022868
$ws.Range("H15").Value = 1323.068604
$ws.Range("I15").Value = 0.4042542905775846
$ws.Range("J15").Value = 0.4112336327632385
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 0.821983
$ws.Range("N15").Value = 2.465949
$ws.Range("O15").Value = 0.03303991048491191
$ws.Range("P15").Value = 0.03750663509291002
$ws.Range("Q15").Value = 362.513300107244
$ws.Range("R15").Value = 3262.619700965196
$ws.Range("S15").Value = 0.01335652557382496
$ws.Range("T15").Value = 0.01542398980198255

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 441.022868
$ws.Range("H16").Value = 1323.068604
$ws.Range("I16").Value = 0.4042542905775846
$ws.Range("J16").Value = 0.4112336327632385
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 7.086962666666667
$ws.Range("N16").Value = 21.260888
$ws.Range("O16").Value = 0.2848630836849173
$ws.Range("P16").Value = 0.3233742335981926
$ws.Range("Q16").Value = 3125.512600662262
$ws.Range("R16").Value = 28129.61340596035
$ws.Range("S16").Value = 0.1151571238067893
$ws.Range("T16").Value = 0.1329823608246128

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 441.022868
$ws.Range("H17").Value = 1323.068604
$ws.Range("I17").Value = 0.4042542905775846
$ws.Range("J17").Value = 0.4112336327632385
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 6.510773
$ws.Range("N17").Value = 19.532319
$ws.Range("O17").Value = 0.2617029270770581
$ws.Range("P17").Value = 0.2970830139841956
$ws.Range("Q17").Value = 2871.399781356964
$ws.Range("R17").Value = 25842.59803221268
$ws.Range("S17").Value = 0.1057945311276134
$ws.Range("T17").Value = 0.1221705270729727

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 441.022868
$ws.Range("H18").Value = 1323.068604
$ws.Range("I18").Value = 0.4042542905775846
$ws.Range("J18").Value = 0.4112336327632385
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 1.570313333333333
$ws.Range("N18").Value = 4.71094
$ws.Range("O18").Value = 0.06311932481158002
$ws.Range("P18").Value = 0.07165253925551318
$ws.Range("Q18").Value = 692.5440899253066
$ws.Range("R18").Value = 6232.89680932776
$ws.Range("S18").Value = 0.02551625787344141
$ws.Range("T18").Value = 0.02946593401475524

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 441.022868
$ws.Range("H19").Value = 1323.068604
$ws.Range("I19").Value = 0.4042542905775846
$ws.Range("J19").Value = 0.4112336327632385
$ws.Range("K19").Value = 2
$ws.Range("M19").Value = 8.3884095
$ws.Range("N19").Value = 16.776819
$ws.Range("O19").Value = 0.3371752201575759
$ws.Range("P19").Value = 0.2551723609258745
$ws.Range("Q19").Value = 3699.480415648446
$ws.Range("R19").Value = 22196.88249389068
$ws.Range("S19").Value = 0.1363045294251417
$ws.Range("T19").Value = 0.1049354569643196

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 536.9901530000001
$ws.Range("H20").Value = 1610.970459
$ws.Range("I20").Value = 0.4922206740267346
$ws.Range("J20").Value = 0.5007187322909461
$ws.Range("K20").Value = 2
$ws.Range("M20").Value = 0.500046
$ws.Range("N20").Value = 1.000092
$ws.Range("O20").Value = 0.02009953378395692
$ws.Range("P20").Value = 0.01521121714331422
$ws.Range("Q20").Value = 268.519778047038
$ws.Range("R20").Value = 1611.118668282228
$ws.Range("S20").Value = 0.009893406066762402
$ws.Range("T20").Value = 0.007616541364602603

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 536.9901530000001
$ws.Range("H21").Value = 1610.970459
$ws.Range("I21").Value = 0.4922206740267346
$ws.Range("J21").Value = 0.5007187322909461
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 0.821983
$ws.Range("N21").Value = 2.465949
$ws.Range("O21").Value = 0.03303991048491191
$ws.Range("P21").Value = 0.03750663509291002
$ws.Range("Q21").Value = 441.3967769333991
$ws.Range("R21").Value = 3972.570992400591
$ws.Range("S21").Value = 0.01626292700866632
$ws.Range("T21").Value = 0.01878027477622101

$ws.Range("E22").Value = 3
$ws.Range("G22").Value = 536.9901530000001
$ws.Range("H22").Value = 1610.970459
$ws.Range("I22").Value = 0.4922206740267346
$ws.Range("J22").Value = 0.5007187322909461
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 7.086962666666667
$ws.Range("N22").Value = 21.260888
$ws.Range("O22").Value = 0.2848630836849173
$ws.Range("P22").Value = 0.3233742335981926
$ws.Range("Q22").Value = 3805.629166678622
$ws.Range("R22").Value = 34250.6625001076
$ws.Range("S22").Value = 0.1402154990567241
$ws.Range("T22").Value = 0.1619195363028433

$ws.Range("E23").Value = 3
$ws.Range("G23").Value = 536.9901530000001
$ws.Range("H23").Value = 1610.970459
$ws.Range("I23").Value = 0.4922206740267346
$ws.Range("J23").Value = 0.5007187322909461
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 6.510773
$ws.Range("N23").Value = 19.532319
$ws.Range("O23").Value = 0.2617029270770581
$ws.Range("P23").Value = 0.2970830139841956
$ws.Range("Q23").Value = 3496.22098941827
$ws.Range("R23").Value = 31465.98890476442
$ws.Range("S23").Value = 0.1288155911606389
$ws.Range("T23").Value = 0.1487550301473398

$ws.Range("E24").Value = 3
$ws.Range("G24").Value = 536.9901530000001
$ws.Range("H24").Value = 1610.970459
$ws.Range("I24").Value = 0.4922206740267346
$ws.Range("J24").Value = 0.5007187322909461
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 1.570313333333333
$ws.Range("N24").Value = 4.71094
$ws.Range("O24").Value = 0.06311932481158002
$ws.Range("P24").Value = 0.07165253925551318
$ws.Range("Q24").Value = 843.2427971246067
$ws.Range("R24").Value = 7589.18517412146
$ws.Range("S24").Value = 0.03106863660286831
$ws.Range("T24").Value = 0.03587776862144781

$ws.Range("E25").Value = 3
$ws.Range("G25").Value = 536.9901530000001
$ws.Range("H25").Value = 1610.970459
$ws.Range("I25").Value = 0.4922206740267346
$ws.Range("J25").Value = 0.5007187322909461
$ws.Range("K25").Value = 2
$ws.Range("M25").Value = 8.3884095
$ws.Range("N25").Value = 16.776819
$ws.Range("O25").Value = 0.3371752201575759
$ws.Range("P25").Value = 0.2551723609258745
$ws.Range("Q25").Value = 4504.493300831654
$ws.Range("R25").Value = 27026.95980498992
$ws.Range("S25").Value = 0.1659646141310747
$ws.Range("T25").Value = 0.1277695810784917

$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 29.47333466666667
$ws.Range("H26").Value = 88.42000400000001
$ws.Range("I26").Value = 0.02701610927946045
$ws.Range("J26").Value = 0.02748253517915091
$ws.Range("K26").Value = 2
$ws.Range("M26").Value = 0.500046
$ws.Range("N26").Value = 1.000092
$ws.Range("O26").Value = 0.02009953378395692
$ws.Range("P26").Value = 0.01521121714331422
$ws.Range("Q26").Value = 14.738023106728
$ws.Range("R26").Value = 88.428138640368
$ws.Range("S26").Value = 0.0005430112011735875
$ws.Range("T26").Value = 0.0004180428102588365

$ws.Range("E27").Value = 3
$ws.Range("G27").Value = 29.47333466666667
$ws.Range("H27").Value = 88.42000400000001
$ws.Range("I27").Value = 0.02701610927946045
$ws.Range("J27").Value = 0.02748253517915091
$ws.Range("K27").Value = 3
$ws.Range("M27").Value = 0.821983
$ws.Range("N27").Value = 2.465949
$ws.Range("O27").Value = 0.03303991048491191
$ws.Range("P27").Value = 0.03750663509291002
$ws.Range("Q27").Value = 24.22658004931067
$ws.Range("R27").Value = 218.039220443796
$ws.Range("S27").Value = 0.0008926098322439716
$ws.Range("T27").Value = 0.001030777418392475

$ws.Range("E28").Value = 3
$ws.Range("G28").Value = 29.47333466666667
$ws.Range("H28").Value = 88.42000400000001
$ws.Range("I28").Value = 0.02701610927946045
$ws.Range("J28").Value = 0.02748253517915091
$ws.Range("K28").Value = 3
$ws.Range("M28").Value = 7.086962666666667
$ws.Range("N28").Value = 21.260888
$ws.Range("O28").Value = 0.2848630836849173
$ws.Range("P28").Value = 0.3233742335981926
$ws.Range("Q28").Value = 208.8764224448392
$ws.Range("R28").Value = 1879.887802003552
$ws.Range("S28").Value = 0.007695892198515813
$ws.Range("T28").Value = 0.00888714375089329

$ws.Range("E29").Value = 3
$ws.Range("G29").Value = 29.47333466666667
$ws.Range("H29").Value = 88.42000400000001
$ws.Range("I29").Value = 0.02701610927946045
$ws.Range("J29").Value = 0.02748253517915091
$ws.Range("K29").Value = 3
$ws.Range("M29").Value = 6.510773
$ws.Range("N29").Value = 19.532319
$ws.Range("O29").Value = 0.2617029270770581
$ws.Range("P29").Value = 0.2970830139841956
$ws.Range("Q29").Value = 191.8941915676974
$ws.Range("R29").Value = 1727.047724109276
$ws.Range("S29").Value = 0.007070194876668471
$ws.Range("T29").Value = 0.008164594382948837

$ws.Range("E30").Value = 3
$ws.Range("G30").Value = 29.47333466666667
$ws.Range("H30").Value = 88.42000400000001
$ws.Range("I30").Value = 0.02701610927946045
$ws.Range("J30").Value = 0.02748253517915091
$ws.Range("K30").Value = 3
$ws.Range("M30").Value = 1.570313333333333
$ws.Range("N30").Value = 4.71094
$ws.Range("O30").Value = 0.06311932481158002
$ws.Range("P30").Value = 0.07165253925551318
$ws.Range("Q30").Value = 46.28237040486223
$ws.Range("R30").Value = 416.54133364376
$ws.Range("S30").Value = 0.001705238576755405
$ws.Range("T30").Value = 0.001969193430765132

$ws.Range("E31").Value = 3
$ws.Range("G31").Value = 29.47333466666667
$ws.Range("H31").Value = 88.42000400000001
$ws.Range("I31").Value = 0.02701610927946045
$ws.Range("J31").Value = 0.02748253517915091
$ws.Range("K31").Value = 2
$ws.Range("M31").Value = 8.3884095
$ws.Range("N31").Value = 16.776819
$ws.Range("O31").Value = 0.3371752201575759
$ws.Range("P31").Value = 0.2551723609258745
$ws.Range("Q31").Value = 247.234400514546
$ws.Range("R31").Value = 1483.406403087276
$ws.Range("S31").Value = 0.009109162594103208
$ws.Range("T31").Value = 0.00701278338589234

$ws.Range("E32").Value = 2
$ws.Range("G32").Value = 13.215904
$ws.Range("H32").Value = 26.431808
$ws.Range("I32").Value = 0.01211407907279189
$ws.Range("J32").Value = 0.008215483604915494
$ws.Range("K32").Value = 2
$ws.Range("M32").Value = 0.500046
$ws.Range("N32").Value = 1.000092
$ws.Range("O32").Value = 0.02009953378395692
$ws.Range("P32").Value = 0.01521121714331422
$ws.Range("Q32").Value = 6.608559931584
$ws.Range("R32").Value = 26.434239726336
$ws.Range("S32").Value = 0.0002434873415851062
$ws.Range("T32").Value = 0.0001249675050517075

$ws.Range("E33").Value = 2
$ws.Range("G33").Value = 13.215904
$ws.Range("H33").Value = 26.431808
$ws.Range("I33").Value = 0.01211407907279189
$ws.Range("J33").Value = 0.008215483604915494
$ws.Range("K33").Value = 3
$ws.Range("M33").Value = 0.821983
$ws.Range("N33").Value = 2.465949
$ws.Range("O33").Value = 0.03303991048491191
$ws.Range("P33").Value = 0.03750663509291002
$ws.Range("Q33").Value = 10.863248417632
$ws.Range("R33").Value = 65.179490505792
$ws.Range("S33").Value = 0.0004002480881721889
$ws.Range("T33").Value = 0.0003081351456813504

$ws.Range("E34").Value = 2
$ws.Range("G34").Value = 13.215904
$ws.Range("H34").Value = 26.431808
$ws.Range("I34").Value = 0.01211407907279189
$ws.Range("J34").Value = 0.008215483604915494
$ws.Range("K34").Value = 3
$ws.Range("M34").Value = 7.086962666666667
$ws.Range("N34").Value = 21.260888
$ws.Range("O34").Value = 0.2848630836849173
$ws.Range("P34").Value = 0.3233742335981926
$ws.Range("Q34").Value = 93.66061825425068
$ws.Range("R34").Value = 561.963709525504
$ws.Range("S34").Value = 0.003450853920678422
$ws.Range("T34").Value = 0.002656675714378064

$ws.Range("E35").Value = 2
$ws.Range("G35").Value = 13.215904
$ws.Range("H35").Value = 26.431808
$ws.Range("I35").Value = 0.01211407907279189
$ws.Range("J35").Value = 0.008215483604915494
$ws.Range("K35").Value = 3
$ws.Range("M35").Value = 6.510773
$ws.Range("N35").Value = 19.532319
$ws.Range("O35").Value = 0.2617029270770581
$ws.Range("P35").Value = 0.2970830139841956
$ws.Range("Q35").Value = 86.04575093379201
$ws.Range("R35").Value = 516.274505602752
$ws.Range("S35").Value = 0.003170289952192572
$ws.Range("T35").Value = 0.002440680630686039

$ws.Range("E36").Value = 2
$ws.Range("G36").Value = 13.215904
$ws.Range("H36").Value = 26.431808
$ws.Range("I36").Value = 0.01211407907279189
$ws.Range("J36").Value = 0.008215483604915494
$ws.Range("K36").Value = 3
$ws.Range("M36").Value = 1.570313333333333
$ws.Range("N36").Value = 4.71094
$ws.Range("O36").Value = 0.06311932481158002
$ws.Range("P36").Value = 0.07165253925551318
$ws.Range("Q36").Value = 20.75311026325333
$ws.Range("R36").Value = 124.51866157952
$ws.Range("S36").Value = 0.0007646324917887157
$ws.Range("T36").Value = 0.0005886602615042323

$ws.Range("E37").Value = 2
$ws.Range("G37").Value = 13.215904
$ws.Range("H37").Value = 26.431808
$ws.Range("I37").Value = 0.01211407907279189
$ws.Range("J37").Value = 0.008215483604915494
$ws.Range("K37").Value = 2
$ws.Range("M37").Value = 8.3884095
$ws.Range("N37").Value = 16.776819
$ws.Range("O37").Value = 0.3371752201575759
$ws.Range("P37").Value = 0.2551723609258745
$ws.Range("Q37").Value = 110.860414664688
$ws.Range("R37").Value = 443.441658658752
$ws.Range("S37").Value = 0.00408456727837489
$ws.Range("T37").Value = 0.002096364347614101
